# Experiment2-Table.xlsx -- "added images and tables, experiment 1-4"
#
# On sheet "Sheet1":
#  - header label in G1 changes from "Missingness 98.5%" text stays the same
#    but the underlying shared-string slot is refreshed (handled automatically
#    by re-assigning the value)
#  - the repeating "smape_mean" / "smape_std" column headers in row 2
#    (C2,E2,G2 and D2,F2,H2) are re-typed with a space instead of an
#    underscore: "smape mean" / "smape std"
#  - the "5-imputations-combined" row labels (B3, B9, B15) are re-typed with
#    spaces instead of hyphens: "5 imputations combined"
#  - the footer note in B26 is refreshed as well
#  - the current selection on the sheet moves from B27 down to the P6:P8
#    block

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("G1").Value = "Missingness 98.5%"

$ws.Range("C2").Value = "smape mean"
$ws.Range("D2").Value = "smape std"
$ws.Range("E2").Value = "smape mean"
$ws.Range("F2").Value = "smape std"
$ws.Range("G2").Value = "smape mean"
$ws.Range("H2").Value = "smape std"

$ws.Range("B3").Value = "5 imputations combined"
$ws.Range("B9").Value = "5 imputations combined"
$ws.Range("B15").Value = "5 imputations combined"

$ws.Range("B26").Value = "WITH SHIFTED  DATA AND DIFFERENT MISSINGNESS"

# Move the sheet's active selection from B27 to P6:P8
$ws.Range("P6:P8").Select()
